$wb = $excel.ActiveWorkbook

# --- 1. Belgium sheet: the data range A1:D15 gets selected (was a whole-sheet
#        selection before). We do this before touching Portugal so Belgium does
#        not end up being the final active sheet.
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Activate()
[void]$belgium.Range("A1:D15").Select()

# --- 2. Build the new "Portugal" sheet as a copy of "Swiss" (keeps styles,
#        merged cells, borders, etc. identical to the other market sheets),
#        placed immediately after it.
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy([System.Reflection.Missing]::Value, $swiss)
$portugal = $wb.Worksheets.Item($swiss.Index + 1)
$portugal.Name = "Portugal"

# --- 3. Market-specific content for Portugal.
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2410"

# --- 4. Column widths specific to the Portugal sheet.
$portugal.Columns.Item(1).ColumnWidth = 26.166666666666668
$portugal.Columns.Item(2).ColumnWidth = 19.166666666666668
$portugal.Columns.Item(3).ColumnWidth = 18.498697916666668
$portugal.Columns.Item(4).ColumnWidth = 13.944010416666666

# --- 5. Row heights: rows 3-5 grow (wrapped ticket/user-story text), row 15
#        goes back to the sheet's default (no explicit custom height).
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8
$portugal.Rows.Item(15).AutoFit()

# --- 6. Final selection/active state: Portugal is the active sheet with B4
#        selected (matches the last-saved UI state captured in the workbook).
$portugal.Activate()
[void]$portugal.Range("B4").Select()
